$wb = $excel.ActiveWorkbook

$ALC = $wb.Worksheets.Item("ALC")
$ARM = $wb.Worksheets.Item("ARM")
$BSM = $wb.Worksheets.Item("BSM")
$CRP = $wb.Worksheets.Item("CRP")
$CUL = $wb.Worksheets.Item("CUL")
$GSM = $wb.Worksheets.Item("GSM")
$LTW = $wb.Worksheets.Item("LTW")
$WVR = $wb.Worksheets.Item("WVR")

# Row 19 (ALC)
$ALC.Range("H19").Value2 = 810.9231
$ALC.Range("I19").Value2 = 652.2222
$ALC.Range("J19").Value2 = 1168
$ALC.Range("K19").Value2 = 652.2222
$ALC.Range("L19").Value2 = 1168
$ALC.Range("M19").Value2 = -477.2222
$ALC.Range("N19").Value2 = -1518

# Row 43 (ALC)
$ALC.Range("H43").Value2 = 15598
$ALC.Range("I43").Value2 = 19649.334
$ALC.Range("J43").Value2 = 3444
$ALC.Range("K43").Value2 = 19649.334
$ALC.Range("L43").Value2 = 3444
$ALC.Range("M43").Value2 = -19580.334
$ALC.Range("N43").Value2 = -3582

# Row 53 (ALC)
$ALC.Range("H53").Value2 = 218.21428
$ALC.Range("I53").Value2 = 304.25
$ALC.Range("J53").Value2 = 183.8
$ALC.Range("K53").Value2 = 304.25
$ALC.Range("L53").Value2 = 183.8
$ALC.Range("M53").Value2 = 332.75
$ALC.Range("N53").Value2 = -1457.8

# Row 97 (ALC)
$ALC.Range("H97").Value2 = 1690.7142
$ALC.Range("I97").Value2 = 0
$ALC.Range("J97").Value2 = 1690.7142
$ALC.Range("K97").Value2 = 0
$ALC.Range("L97").Value2 = 5072.142599999999
$ALC.Range("N97").Value2 = -6064.142599999999

# Row 141 (ALC)
$ALC.Range("H141").Value2 = 2370.1333
$ALC.Range("I141").Value2 = 1713
$ALC.Range("J141").Value2 = 4998.6665
$ALC.Range("K141").Value2 = 5139
$ALC.Range("L141").Value2 = 14995.9995
$ALC.Range("M141").Value2 = 41

# Row 5 (ARM)
$ARM.Range("H5").Value2 = 2341.8096
$ARM.Range("I5").Value2 = 255.71428
$ARM.Range("J5").Value2 = 6514
$ARM.Range("K5").Value2 = 255.71428
$ARM.Range("L5").Value2 = 6514
$ARM.Range("M5").Value2 = -143.71428

# Row 61 (ARM)
$ARM.Range("H61").Value2 = 3034.2222
$ARM.Range("I61").Value2 = 2763.4
$ARM.Range("J61").Value2 = 5200.8
$ARM.Range("K61").Value2 = 2763.4
$ARM.Range("L61").Value2 = 5200.8
$ARM.Range("M61").Value2 = -2551.4

# Row 122 (ARM)
$ARM.Range("H122").Value2 = 2306.7778
$ARM.Range("I122").Value2 = 1924.6666
$ARM.Range("J122").Value2 = 3071
$ARM.Range("K122").Value2 = 5773.9998
$ARM.Range("L122").Value2 = 9213
$ARM.Range("M122").Value2 = -3323.9998

# Row 132 (ARM)
$ARM.Range("H132").Value2 = 432329.78
$ARM.Range("I132").Value2 = 464167.78
$ARM.Range("J132").Value2 = 2517
$ARM.Range("K132").Value2 = 1392503.34
$ARM.Range("L132").Value2 = 7551
$ARM.Range("M132").Value2 = -1389973.34
$ARM.Range("N132").Value2 = -12611

# Row 136 (ARM)
$ARM.Range("H136").Value2 = 3034.2222
$ARM.Range("I136").Value2 = 2763.4
$ARM.Range("J136").Value2 = 5200.8
$ARM.Range("K136").Value2 = 8290.200000000001
$ARM.Range("L136").Value2 = 15602.4
$ARM.Range("M136").Value2 = -5740.200000000001

# Row 4 (BSM)
$BSM.Range("H4").Value2 = 2341.8096
$BSM.Range("I4").Value2 = 255.71428
$BSM.Range("J4").Value2 = 6514
$BSM.Range("K4").Value2 = 255.71428
$BSM.Range("L4").Value2 = 6514
$BSM.Range("M4").Value2 = -140.71428

# Row 134 (BSM)
$BSM.Range("H134").Value2 = 4937.375
$BSM.Range("I134").Value2 = 2397.0967
$BSM.Range("J134").Value2 = 13687.223
$BSM.Range("K134").Value2 = 7191.2901
$BSM.Range("L134").Value2 = 41061.669
$BSM.Range("M134").Value2 = -4656.2901
$BSM.Range("N134").Value2 = -46131.669

# Row 9 (CRP)
$CRP.Range("H9").Value2 = 21808.166
$CRP.Range("I9").Value2 = 0
$CRP.Range("J9").Value2 = 21808.166
$CRP.Range("K9").Value2 = 0
$CRP.Range("L9").Value2 = 21808.166
$CRP.Range("N9").Value2 = -22144.166

# Row 10 (CRP)
$CRP.Range("H10").Value2 = 1373.1177
$CRP.Range("I10").Value2 = 1290.8889
$CRP.Range("J10").Value2 = 1465.625
$CRP.Range("K10").Value2 = 1290.8889
$CRP.Range("L10").Value2 = 1465.625
$CRP.Range("M10").Value2 = -1151.8889
$CRP.Range("N10").Value2 = -1743.625

# Row 31 (CRP)
$CRP.Range("H31").Value2 = 2497.577
$CRP.Range("I31").Value2 = 2382.0952
$CRP.Range("J31").Value2 = 2982.6
$CRP.Range("K31").Value2 = 2382.0952
$CRP.Range("L31").Value2 = 2982.6
$CRP.Range("M31").Value2 = -2087.0952
$CRP.Range("N31").Value2 = -3572.6

# Row 34 (CRP)
$CRP.Range("H34").Value2 = 2497.577
$CRP.Range("I34").Value2 = 2382.0952
$CRP.Range("J34").Value2 = 2982.6
$CRP.Range("K34").Value2 = 2382.0952
$CRP.Range("L34").Value2 = 2982.6
$CRP.Range("M34").Value2 = -2180.0952
$CRP.Range("N34").Value2 = -3386.6

# Row 58 (CRP)
$CRP.Range("H58").Value2 = 9741.093999999999
$CRP.Range("I58").Value2 = 7506.1665
$CRP.Range("J58").Value2 = 10256.846
$CRP.Range("K58").Value2 = 7506.1665
$CRP.Range("L58").Value2 = 10256.846
$CRP.Range("M58").Value2 = -7303.1665
$CRP.Range("N58").Value2 = -10662.846

# Row 86 (CRP)
$CRP.Range("H86").Value2 = 100000
$CRP.Range("I86").Value2 = 100000
$CRP.Range("J86").Value2 = 0
$CRP.Range("K86").Value2 = 100000
$CRP.Range("L86").Value2 = 0
$CRP.Range("M86").Value2 = -98877
$CRP.Range("N86").ClearContents()

# Row 89 (CRP)
$CRP.Range("H89").Value2 = 100000
$CRP.Range("I89").Value2 = 100000
$CRP.Range("J89").Value2 = 0
$CRP.Range("K89").Value2 = 500000
$CRP.Range("L89").Value2 = 0
$CRP.Range("M89").Value2 = -494384
$CRP.Range("N89").ClearContents()

# Row 99 (CRP)
$CRP.Range("H99").Value2 = 3384.625
$CRP.Range("I99").Value2 = 3144.25
$CRP.Range("J99").Value2 = 3625
$CRP.Range("K99").Value2 = 3144.25
$CRP.Range("L99").Value2 = 3625
$CRP.Range("M99").Value2 = -1646.25
$CRP.Range("N99").Value2 = -6621

# Row 122 (CRP)
$CRP.Range("H122").Value2 = 13109
$CRP.Range("I122").Value2 = 1635.7778
$CRP.Range("J122").Value2 = 30318.834
$CRP.Range("K122").Value2 = 4907.3334
$CRP.Range("L122").Value2 = 90956.50199999999
$CRP.Range("M122").Value2 = -2457.3334

# Row 126 (CRP)
$CRP.Range("H126").Value2 = 3384.625
$CRP.Range("I126").Value2 = 3144.25
$CRP.Range("J126").Value2 = 3625
$CRP.Range("K126").Value2 = 9432.75
$CRP.Range("L126").Value2 = 10875
$CRP.Range("M126").Value2 = -6962.75
$CRP.Range("N126").Value2 = -15815

# Row 132 (CRP)
$CRP.Range("H132").Value2 = 3327.8
$CRP.Range("I132").Value2 = 3221.375
$CRP.Range("J132").Value2 = 3753.5
$CRP.Range("K132").Value2 = 9664.125
$CRP.Range("L132").Value2 = 11260.5
$CRP.Range("M132").Value2 = -7134.125

# Row 134 (CRP)
$CRP.Range("H134").Value2 = 1973.5
$CRP.Range("I134").Value2 = 1507.8572
$CRP.Range("J134").Value2 = 5233
$CRP.Range("K134").Value2 = 4523.571599999999
$CRP.Range("L134").Value2 = 15699
$CRP.Range("M134").Value2 = -1988.571599999999

# Row 136 (CRP)
$CRP.Range("H136").Value2 = 9741.093999999999
$CRP.Range("I136").Value2 = 7506.1665
$CRP.Range("J136").Value2 = 10256.846
$CRP.Range("K136").Value2 = 22518.4995
$CRP.Range("L136").Value2 = 30770.538
$CRP.Range("M136").Value2 = -19968.4995
$CRP.Range("N136").Value2 = -35870.538

# Row 68 (CUL)
$CUL.Range("H68").Value2 = 4926.357
$CUL.Range("I68").Value2 = 0
$CUL.Range("J68").Value2 = 4926.357
$CUL.Range("K68").Value2 = 0
$CUL.Range("L68").Value2 = 14779.071
$CUL.Range("N68").Value2 = -16401.071

# Row 71 (CUL)
$CUL.Range("H71").Value2 = 4926.357
$CUL.Range("I71").Value2 = 0
$CUL.Range("J71").Value2 = 4926.357
$CUL.Range("K71").Value2 = 0
$CUL.Range("L71").Value2 = 44337.213
$CUL.Range("N71").Value2 = -52449.213

# Row 18 (GSM)
$GSM.Range("H18").Value2 = 9965.333000000001
$GSM.Range("I18").Value2 = 9998
$GSM.Range("J18").Value2 = 9900
$GSM.Range("K18").Value2 = 9998
$GSM.Range("L18").Value2 = 9900
$GSM.Range("M18").Value2 = -9705

# Row 127 (GSM)
$GSM.Range("H127").Value2 = 19975
$GSM.Range("I127").Value2 = 0
$GSM.Range("J127").Value2 = 19975
$GSM.Range("K127").Value2 = 0
$GSM.Range("L127").Value2 = 19975
$GSM.Range("N127").Value2 = -29895

# Row 132 (GSM)
$GSM.Range("H132").Value2 = 11042.625
$GSM.Range("I132").Value2 = 11616.486
$GSM.Range("J132").Value2 = 3965
$GSM.Range("K132").Value2 = 34849.458
$GSM.Range("L132").Value2 = 11895
$GSM.Range("M132").Value2 = -32319.458

# Row 7 (LTW)
$LTW.Range("H7").Value2 = 0
$LTW.Range("I7").Value2 = 0
$LTW.Range("J7").Value2 = 0
$LTW.Range("K7").Value2 = 0
$LTW.Range("L7").Value2 = 0
$LTW.Range("M7").ClearContents()

# Row 40 (LTW)
$LTW.Range("H40").Value2 = 4206.75
$LTW.Range("I40").Value2 = 3666
$LTW.Range("J40").Value2 = 4747.5
$LTW.Range("K40").Value2 = 3666
$LTW.Range("L40").Value2 = 4747.5
$LTW.Range("M40").Value2 = -3530

# Row 93 (LTW)
$LTW.Range("H93").Value2 = 2382.2727
$LTW.Range("I93").Value2 = 1862.7142
$LTW.Range("J93").Value2 = 3291.5
$LTW.Range("K93").Value2 = 1862.7142
$LTW.Range("L93").Value2 = 3291.5
$LTW.Range("M93").Value2 = -614.7141999999999
$LTW.Range("N93").Value2 = -5787.5

# Row 122 (LTW)
$LTW.Range("H122").Value2 = 7799
$LTW.Range("I122").Value2 = 5000
$LTW.Range("J122").Value2 = 8498.75
$LTW.Range("K122").Value2 = 15000
$LTW.Range("L122").Value2 = 25496.25
$LTW.Range("M122").Value2 = -12550

# Row 126 (LTW)
$LTW.Range("H126").Value2 = 0
$LTW.Range("I126").Value2 = 0
$LTW.Range("J126").Value2 = 0
$LTW.Range("K126").Value2 = 0
$LTW.Range("L126").Value2 = 0
$LTW.Range("M126").ClearContents()

# Row 126 (WVR)
$WVR.Range("H126").Value2 = 2257.0908
$WVR.Range("I126").Value2 = 2257.0908
$WVR.Range("J126").Value2 = 0
$WVR.Range("K126").Value2 = 6771.2724
$WVR.Range("L126").Value2 = 0
$WVR.Range("M126").Value2 = -4301.2724
